$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder columns D:G -------------------------------------------------
# Current layout:  D=Notes        E=MaxBookings  F=Booked  G=Available
# Target layout:   D=MaxBookings  E=Booked       F=Available  G=Notes
# i.e. shift MaxBookings/Booked/Available one column left, and move the
# Notes column to the end (G), for the header row and every data row.

$lastRow = 9

for ($r = 1; $r -le $lastRow; $r++) {
    $notes = $ws.Cells.Item($r, 4).Value()
    $maxBookings = $ws.Cells.Item($r, 5).Value()
    $booked = $ws.Cells.Item($r, 6).Value()
    $available = $ws.Cells.Item($r, 7).Value()

    $ws.Cells.Item($r, 4).Value = $maxBookings
    $ws.Cells.Item($r, 5).Value = $booked
    $ws.Cells.Item($r, 6).Value = $available
    $ws.Cells.Item($r, 7).Value = $notes
}

# --- Append two new rows for the newly synced calendar days --------------
# Row 10: 2025-02-16 (serial 45704) - Closed / Maintenance
# Row 11: 2025-02-17 (serial 45705) - Closed / Maintenance

$ws.Cells.Item(10, 1).Value = 45704
$ws.Cells.Item(10, 2).Value = "Closed"
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = "Maintenance"

$ws.Cells.Item(11, 1).Value = 45705
$ws.Cells.Item(11, 2).Value = "Closed"
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = "Maintenance"

# Match the source date formatting/style used by the other date cells (A2:A9)
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range("A10:A11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("D16").Select() | Out-Null
